$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tblIngredients")

# Add the new ingredient row (48): "grated parmesan", type "Grocery", all nutrition values 0
$ws.Range("A48").Value = "grated parmesan"
$ws.Range("B48").Value = "Grocery"
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0

# Update the hidden _FilterDatabase defined name to cover the new row
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -eq "tblIngredients!_FilterDatabase") {
        $n.RefersTo = "=tblIngredients!`$A`$1:`$G`$48"
    }
}

# Turn on the AutoFilter over the full new range first so the autoFilter
# element picks up A1:G48 as its reference range...
$ws.Range("A1:G48").AutoFilter()

# ...then apply the actual column filter, showing only "Grocery" rows
$ws.Range("A1:G48").AutoFilter(2, @("Grocery"))

# Update the active selection shown in the sheet view
$ws.Activate()
$ws.Range("H51").Select()
